# Auto-generated edit script: update cryptos price/volume data and fix a row-43/44 ordering swap
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write each updated value as a literal-string formula so Excel treats it as text
# (avoids turning numeric-looking strings like "7.74" into real numbers).
$ws.Range("D2").Formula = "=`"42.884.69`""
$ws.Range("E2").Formula = "=`"  -2.00%  `""
$ws.Range("D3").Formula = "=`"2.561.06`""
$ws.Range("E3").Formula = "=`"  -3.22%  `""
$ws.Range("E4").Formula = "=`"  -0.01%  `""
$ws.Range("D5").Formula = "=`"301.12`""
$ws.Range("E5").Formula = "=`"  -0.94%  `""
$ws.Range("D6").Formula = "=`"94.39`""
$ws.Range("E6").Formula = "=`"  -2.62%  `""
$ws.Range("D7").Formula = "=`"0.572`""
$ws.Range("E7").Formula = "=`"  -2.57%  `""
$ws.Range("E8").Formula = "=`"  +0.02%  `""
$ws.Range("E9").Formula = "=`"  -3.19%  `""
$ws.Range("D10").Formula = "=`"36.32`""
$ws.Range("E10").Formula = "=`"  -3.09%  `""
$ws.Range("D11").Formula = "=`"0.0810`""
$ws.Range("E11").Formula = "=`"  -1.18%  `""
$ws.Range("D12").Formula = "=`"7.74`""
$ws.Range("E12").Formula = "=`"  -1.90%  `""
$ws.Range("D13").Formula = "=`"0.114`""
$ws.Range("E13").Formula = "=`"  +6.91%  `""
$ws.Range("D14").Formula = "=`"2.955.72`""
$ws.Range("E14").Formula = "=`"  -3.65%  `""
$ws.Range("D15").Formula = "=`"2.592.83`""
$ws.Range("E15").Formula = "=`"  -2.34%  `""
$ws.Range("D16").Formula = "=`"0.878`""
$ws.Range("E16").Formula = "=`"  -2.39%  `""
$ws.Range("E17").Formula = "=`"  -3.35%  `""
$ws.Range("D18").Formula = "=`"42.924.28`""
$ws.Range("E18").Formula = "=`"  -1.81%  `""
$ws.Range("D19").Formula = "=`"0.0₃0986`""
$ws.Range("E19").Formula = "=`"  -0.21%  `""
$ws.Range("D20").Formula = "=`"12.68`""
$ws.Range("E20").Formula = "=`"  +0.65%  `""
$ws.Range("E21").Formula = "=`"  -2.98%  `""
$ws.Range("D22").Formula = "=`"71.58`""
$ws.Range("E22").Formula = "=`"  -4.21%  `""
$ws.Range("D23").Formula = "=`"253.00`""
$ws.Range("E23").Formula = "=`"  -7.19%  `""
$ws.Range("D24").Formula = "=`"2.95`""
$ws.Range("E24").Formula = "=`"  -1.05%  `""
$ws.Range("E25").Formula = "=`"  -6.13%  `""
$ws.Range("D26").Formula = "=`"28.89`""
$ws.Range("E26").Formula = "=`"  -4.59%  `""
$ws.Range("E28").Formula = "=`"  -1.68%  `""
$ws.Range("D29").Formula = "=`"37.00`""
$ws.Range("E29").Formula = "=`"  -2.58%  `""
$ws.Range("E30").Formula = "=`"  -4.19%  `""
$ws.Range("D31").Formula = "=`"6.08`""
$ws.Range("E31").Formula = "=`"  -0.35%  `""
$ws.Range("D32").Formula = "=`"154.46`""
$ws.Range("E32").Formula = "=`"  +0.35%  `""
$ws.Range("D33").Formula = "=`"2.76`""
$ws.Range("E33").Formula = "=`"  -1.41%  `""
$ws.Range("D34").Formula = "=`"3.38`""
$ws.Range("E34").Formula = "=`"  -7.80%  `""
$ws.Range("E35").Formula = "=`"  -6.38%  `""
$ws.Range("D36").Formula = "=`"0.0798`""
$ws.Range("E36").Formula = "=`"  -2.90%  `""
$ws.Range("D37").Formula = "=`"0.114`""
$ws.Range("E37").Formula = "=`"  -5.05%  `""
$ws.Range("D38").Formula = "=`"17.75`""
$ws.Range("E38").Formula = "=`"  +10.55%  `""
$ws.Range("D39").Formula = "=`"0.119`""
$ws.Range("E39").Formula = "=`"  -1.80%  `""
$ws.Range("D40").Formula = "=`"23.02`""
$ws.Range("E40").Formula = "=`"  -8.49%  `""
$ws.Range("D41").Formula = "=`"2.14`""
$ws.Range("E41").Formula = "=`"  +33.20%  `""
$ws.Range("D42").Formula = "=`"3.41`""
$ws.Range("E42").Formula = "=`"  -4.14%  `""
$ws.Range("B43").Formula = "=`"VeChain`""
$ws.Range("C43").Formula = "=`"https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet`""
$ws.Range("D43").Formula = "=`"0.0310`""
$ws.Range("E43").Formula = "=`"  -2.91%  `""
$ws.Range("B44").Formula = "=`"RenderToken`""
$ws.Range("C44").Formula = "=`"https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr`""
$ws.Range("D44").Formula = "=`"3.87`""
$ws.Range("E44").Formula = "=`"  -0.03%  `""
$ws.Range("E45").Formula = "=`"  +0.15%  `""
$ws.Range("E46").Formula = "=`"  +0.33%  `""
$ws.Range("E47").Formula = "=`"  +1.18%  `""
$ws.Range("D48").Formula = "=`"85.37`""
$ws.Range("E48").Formula = "=`"  -4.98%  `""
$ws.Range("D49").Formula = "=`"106.01`""
$ws.Range("E49").Formula = "=`"  -2.19%  `""
$ws.Range("D50").Formula = "=`"2.810.79`""
$ws.Range("E50").Formula = "=`"  -3.64%  `""
$ws.Range("D51").Formula = "=`"74.07`""
$ws.Range("E51").Formula = "=`"  +4.69%  `""

# Step 2: flatten the formulas back down to plain static text values by copy-pasting
# each touched cell onto itself (keeps cell type as text, adds no new number formats).
$ws.Range("D2").Copy($ws.Range("D2"))
$ws.Range("E2").Copy($ws.Range("E2"))
$ws.Range("D3").Copy($ws.Range("D3"))
$ws.Range("E3").Copy($ws.Range("E3"))
$ws.Range("E4").Copy($ws.Range("E4"))
$ws.Range("D5").Copy($ws.Range("D5"))
$ws.Range("E5").Copy($ws.Range("E5"))
$ws.Range("D6").Copy($ws.Range("D6"))
$ws.Range("E6").Copy($ws.Range("E6"))
$ws.Range("D7").Copy($ws.Range("D7"))
$ws.Range("E7").Copy($ws.Range("E7"))
$ws.Range("E8").Copy($ws.Range("E8"))
$ws.Range("E9").Copy($ws.Range("E9"))
$ws.Range("D10").Copy($ws.Range("D10"))
$ws.Range("E10").Copy($ws.Range("E10"))
$ws.Range("D11").Copy($ws.Range("D11"))
$ws.Range("E11").Copy($ws.Range("E11"))
$ws.Range("D12").Copy($ws.Range("D12"))
$ws.Range("E12").Copy($ws.Range("E12"))
$ws.Range("D13").Copy($ws.Range("D13"))
$ws.Range("E13").Copy($ws.Range("E13"))
$ws.Range("D14").Copy($ws.Range("D14"))
$ws.Range("E14").Copy($ws.Range("E14"))
$ws.Range("D15").Copy($ws.Range("D15"))
$ws.Range("E15").Copy($ws.Range("E15"))
$ws.Range("D16").Copy($ws.Range("D16"))
$ws.Range("E16").Copy($ws.Range("E16"))
$ws.Range("E17").Copy($ws.Range("E17"))
$ws.Range("D18").Copy($ws.Range("D18"))
$ws.Range("E18").Copy($ws.Range("E18"))
$ws.Range("D19").Copy($ws.Range("D19"))
$ws.Range("E19").Copy($ws.Range("E19"))
$ws.Range("D20").Copy($ws.Range("D20"))
$ws.Range("E20").Copy($ws.Range("E20"))
$ws.Range("E21").Copy($ws.Range("E21"))
$ws.Range("D22").Copy($ws.Range("D22"))
$ws.Range("E22").Copy($ws.Range("E22"))
$ws.Range("D23").Copy($ws.Range("D23"))
$ws.Range("E23").Copy($ws.Range("E23"))
$ws.Range("D24").Copy($ws.Range("D24"))
$ws.Range("E24").Copy($ws.Range("E24"))
$ws.Range("E25").Copy($ws.Range("E25"))
$ws.Range("D26").Copy($ws.Range("D26"))
$ws.Range("E26").Copy($ws.Range("E26"))
$ws.Range("E28").Copy($ws.Range("E28"))
$ws.Range("D29").Copy($ws.Range("D29"))
$ws.Range("E29").Copy($ws.Range("E29"))
$ws.Range("E30").Copy($ws.Range("E30"))
$ws.Range("D31").Copy($ws.Range("D31"))
$ws.Range("E31").Copy($ws.Range("E31"))
$ws.Range("D32").Copy($ws.Range("D32"))
$ws.Range("E32").Copy($ws.Range("E32"))
$ws.Range("D33").Copy($ws.Range("D33"))
$ws.Range("E33").Copy($ws.Range("E33"))
$ws.Range("D34").Copy($ws.Range("D34"))
$ws.Range("E34").Copy($ws.Range("E34"))
$ws.Range("E35").Copy($ws.Range("E35"))
$ws.Range("D36").Copy($ws.Range("D36"))
$ws.Range("E36").Copy($ws.Range("E36"))
$ws.Range("D37").Copy($ws.Range("D37"))
$ws.Range("E37").Copy($ws.Range("E37"))
$ws.Range("D38").Copy($ws.Range("D38"))
$ws.Range("E38").Copy($ws.Range("E38"))
$ws.Range("D39").Copy($ws.Range("D39"))
$ws.Range("E39").Copy($ws.Range("E39"))
$ws.Range("D40").Copy($ws.Range("D40"))
$ws.Range("E40").Copy($ws.Range("E40"))
$ws.Range("D41").Copy($ws.Range("D41"))
$ws.Range("E41").Copy($ws.Range("E41"))
$ws.Range("D42").Copy($ws.Range("D42"))
$ws.Range("E42").Copy($ws.Range("E42"))
$ws.Range("B43").Copy($ws.Range("B43"))
$ws.Range("C43").Copy($ws.Range("C43"))
$ws.Range("D43").Copy($ws.Range("D43"))
$ws.Range("E43").Copy($ws.Range("E43"))
$ws.Range("B44").Copy($ws.Range("B44"))
$ws.Range("C44").Copy($ws.Range("C44"))
$ws.Range("D44").Copy($ws.Range("D44"))
$ws.Range("E44").Copy($ws.Range("E44"))
$ws.Range("E45").Copy($ws.Range("E45"))
$ws.Range("E46").Copy($ws.Range("E46"))
$ws.Range("E47").Copy($ws.Range("E47"))
$ws.Range("D48").Copy($ws.Range("D48"))
$ws.Range("E48").Copy($ws.Range("E48"))
$ws.Range("D49").Copy($ws.Range("D49"))
$ws.Range("E49").Copy($ws.Range("E49"))
$ws.Range("D50").Copy($ws.Range("D50"))
$ws.Range("E50").Copy($ws.Range("E50"))
$ws.Range("D51").Copy($ws.Range("D51"))
$ws.Range("E51").Copy($ws.Range("E51"))
